# Update the correct week
# The "Forecast Comparison" sheet had its Week labels off by 4 weeks
# (it should start at W05 instead of W1). Shift labels and update the
# Seasonality Index (and a couple MyForecast) values that correspond to
# the corrected weeks. Also refresh the two dependent summary figures.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Row -> (Week label, Seasonality Index, optional MyForecast override)
$updates = @(
    @{ Row = 2;  Week = "W05"; Seasonality = 1.01 },
    @{ Row = 3;  Week = "W06"; Seasonality = 1.2 },
    @{ Row = 4;  Week = "W07"; Seasonality = 0.91 },
    @{ Row = 5;  Week = "W08"; Seasonality = 1.08 },
    @{ Row = 6;  Week = "W09"; Seasonality = 1.12 },
    @{ Row = 7;  Week = "W10"; Seasonality = 1.02 },
    @{ Row = 8;  Week = "W11"; Seasonality = 1.03; MyForecast = 69 },
    @{ Row = 9;  Week = "W12"; Seasonality = 0.96 },
    @{ Row = 10; Week = "W13"; Seasonality = 1.17 },
    @{ Row = 11; Week = "W14"; Seasonality = 0.8100000000000001 },
    @{ Row = 12; Week = "W15"; Seasonality = 0.89; MyForecast = 83 },
    @{ Row = 13; Week = "W16"; Seasonality = 0.99; MyForecast = 69 },
    @{ Row = 14; Week = "W17"; Seasonality = 1.06 },
    @{ Row = 15; Week = "W18" },
    @{ Row = 16; Week = "W19"; Seasonality = 0.82 },
    @{ Row = 17; Week = "W20"; Seasonality = 0.9399999999999999; MyForecast = 80 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $wsForecast.Cells.Item($r, 1).Value = $u.Week
    if ($u.ContainsKey("MyForecast")) {
        $wsForecast.Cells.Item($r, 4).Value = $u.MyForecast
    }
    if ($u.ContainsKey("Seasonality")) {
        $wsForecast.Cells.Item($r, 16).Value = $u.Seasonality
    }
}

# Summary sheet dependent values (leading apostrophe keeps these as text,
# matching the original inline-string cell type rather than becoming numbers)
$wsSummary.Range("B10").Value = "'849"
$wsSummary.Range("B14").Value = "'69"
